# Copilot-ROI-Rechner.xlsx — update training budget default from €800 to €1.000
# and refresh the associated hint texts / formulas that reference it.

$wb = $excel.ActiveWorkbook

# ---- Sheet "1. Eingaben" ----------------------------------------------
$ws1 = $wb.Worksheets.Item("1. Eingaben")

# Training-Budget pro Mitarbeiter: 800 -> 1000
$ws1.Range("B23").Value = 1000

# Hinweistext zum Training-Budget
$ws1.Range("D23").Value = "Empfohlen: €800-1.500 (Standard: €1.000)"

# Hinweis zur Adoption-Rate
$ws1.Range("A29").Value = "Mit professionellem Training (ab €1.000/MA) erreichen Sie 60-85% Adoption."

# Hinweis zur Zeitersparnis (Forrester TEI)
$ws1.Range("D33").Value = "Forrester TEI: 9h/Monat (≈ 5,2% bei 40h/Woche)"

# ---- Sheet "3. Nutzen" -------------------------------------------------
$ws3 = $wb.Worksheets.Item("3. Nutzen")

# Szenario-Label "Standard (€800)" -> "Standard (€1.000)"
$ws3.Range("A31").Value = "Standard (€1.000)"

# Zugehörige Adoption-Formel auf Basis von 1000 statt 800
$ws3.Range("B31").Formula = "=MIN(0.9,MAX(0.05, 0.05 + 0.85 * (1 - EXP(-1000/600))))"
